$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Tgfb1"
$ws.Range("C2").Value = "Itgb6"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 98.91277700000001
$ws.Range("H2").Value = 296.738331
$ws.Range("I2").Value = 0.8120825131376513
$ws.Range("J2").Value = 0.8120825131376513
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.128124
$ws.Range("N2").Value = 0.384372
$ws.Range("O2").Value = 0.3522399658364659
$ws.Range("P2").Value = 0.352239965836466
$ws.Range("Q2").Value = 12.673100640348
$ws.Range("R2").Value = 114.057905763132
$ws.Range("S2").Value = 0.2860479166839977
$ws.Range("T2").Value = 0.2860479166839978

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Tgfb1"
$ws.Range("C3").Value = "Itgb6"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 98.91277700000001
$ws.Range("H3").Value = 296.738331
$ws.Range("I3").Value = 0.8120825131376513
$ws.Range("J3").Value = 0.8120825131376513
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.2109236666666666
$ws.Range("N3").Value = 0.632771
$ws.Range("O3").Value = 0.5798737562109268
$ws.Range("P3").Value = 0.5798737562109268
$ws.Range("Q3").Value = 20.86304560502233
$ws.Range("R3").Value = 187.767410445201
$ws.Range("S3").Value = 0.4709053372463392
$ws.Range("T3").Value = 0.4709053372463392

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Tgfb1"
$ws.Range("C4").Value = "Itgb6"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 98.91277700000001
$ws.Range("H4").Value = 296.738331
$ws.Range("I4").Value = 0.8120825131376513
$ws.Range("J4").Value = 0.8120825131376513
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.024693
$ws.Range("N4").Value = 0.07407900000000001
$ws.Range("O4").Value = 0.06788627795260727
$ws.Range("P4").Value = 0.06788627795260727
$ws.Range("Q4").Value = 2.442453202461
$ws.Range("R4").Value = 21.982078822149
$ws.Range("S4").Value = 0.05512925920731444
$ws.Range("T4").Value = 0.05512925920731444

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Tgfb1"
$ws.Range("C5").Value = "Itgb6"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 17.04862266666667
$ws.Range("H5").Value = 51.14586800000001
$ws.Range("I5").Value = 0.1399706767982279
$ws.Range("J5").Value = 0.1399706767982279
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.128124
$ws.Range("N5").Value = 0.384372
$ws.Range("O5").Value = 0.3522399658364659
$ws.Range("P5").Value = 0.352239965836466
$ws.Range("Q5").Value = 2.184337730544001
$ws.Range("R5").Value = 19.65903957489601
$ws.Range("S5").Value = 0.04930326641351482
$ws.Range("T5").Value = 0.04930326641351483

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Tgfb1"
$ws.Range("C6").Value = "Itgb6"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 17.04862266666667
$ws.Range("H6").Value = 51.14586800000001
$ws.Range("I6").Value = 0.1399706767982279
$ws.Range("J6").Value = 0.1399706767982279
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2109236666666666
$ws.Range("N6").Value = 0.632771
$ws.Range("O6").Value = 0.5798737562109268
$ws.Range("P6").Value = 0.5798737562109268
$ws.Range("Q6").Value = 3.595958004469778
$ws.Range("R6").Value = 32.363622040228
$ws.Range("S6").Value = 0.08116532211437406
$ws.Range("T6").Value = 0.08116532211437406

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Tgfb1"
$ws.Range("C7").Value = "Itgb6"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 17.04862266666667
$ws.Range("H7").Value = 51.14586800000001
$ws.Range("I7").Value = 0.1399706767982279
$ws.Range("J7").Value = 0.1399706767982279
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.024693
$ws.Range("N7").Value = 0.07407900000000001
$ws.Range("O7").Value = 0.06788627795260727
$ws.Range("P7").Value = 0.06788627795260727
$ws.Range("Q7").Value = 0.4209816395080002
$ws.Range("R7").Value = 3.788834755572001
$ws.Range("S7").Value = 0.009502088270339059
$ws.Range("T7").Value = 0.009502088270339059

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Tgfb1"
$ws.Range("C8").Value = "Itgb6"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 5.839988000000001
$ws.Range("H8").Value = 17.519964
$ws.Range("I8").Value = 0.0479468100641207
$ws.Range("J8").Value = 0.04794681006412069
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.128124
$ws.Range("N8").Value = 0.384372
$ws.Range("O8").Value = 0.3522399658364659
$ws.Range("P8").Value = 0.352239965836466
$ws.Range("Q8").Value = 0.7482426225120002
$ws.Range("R8").Value = 6.734183602608002
$ws.Range("S8").Value = 0.0168887827389534
$ws.Range("T8").Value = 0.0168887827389534

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Tgfb1"
$ws.Range("C9").Value = "Itgb6"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 5.839988000000001
$ws.Range("H9").Value = 17.519964
$ws.Range("I9").Value = 0.0479468100641207
$ws.Range("J9").Value = 0.04794681006412069
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.2109236666666666
$ws.Range("N9").Value = 0.632771
$ws.Range("O9").Value = 0.5798737562109268
$ws.Range("P9").Value = 0.5798737562109268
$ws.Range("Q9").Value = 1.231791682249333
$ws.Range("R9").Value = 11.086125140244
$ws.Range("S9").Value = 0.02780309685021354
$ws.Range("T9").Value = 0.02780309685021354

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Tgfb1"
$ws.Range("C10").Value = "Itgb6"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 5.839988000000001
$ws.Range("H10").Value = 17.519964
$ws.Range("I10").Value = 0.0479468100641207
$ws.Range("J10").Value = 0.04794681006412069
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.024693
$ws.Range("N10").Value = 0.07407900000000001
$ws.Range("O10").Value = 0.06788627795260727
$ws.Range("P10").Value = 0.06788627795260727
$ws.Range("Q10").Value = 0.144206823684
$ws.Range("R10").Value = 1.297861413156
$ws.Range("S10").Value = 0.003254930474953765
$ws.Range("T10").Value = 0.003254930474953765
